# Updated cryptos list values (price/volume columns D and E, rows 2-51).
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the original inlineStr cells) instead of auto-
# converting numeric-looking strings (e.g. "533.64") into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''58.011.07'
$ws.Range("E2").Value = '''  -2.95%  '
$ws.Range("D3").Value = '''2.572.74'
$ws.Range("E3").Value = '''  -2.35%  '
$ws.Range("E4").Value = '''  -0.04%  '
$ws.Range("D5").Value = '''533.64'
$ws.Range("E5").Value = '''  -0.74%  '
$ws.Range("D6").Value = '''141.04'
$ws.Range("E6").Value = '''  -2.06%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '''  +0.06%  '
$ws.Range("D8").Value = '''0.582'
$ws.Range("E8").Value = '''  +2.69%  '
$ws.Range("E9").Value = '''  +3.08%  '
$ws.Range("E10").Value = '''  -4.38%  '
$ws.Range("E11").Value = '''  +3.06%  '
$ws.Range("E12").Value = '''  -1.90%  '
$ws.Range("D13").Value = '''3.028.03'
$ws.Range("E13").Value = '''  -2.48%  '
$ws.Range("D14").Value = '''57.949.59'
$ws.Range("E14").Value = '''  -2.93%  '
$ws.Range("D15").Value = '''20.61'
$ws.Range("E15").Value = '''  -1.49%  '
$ws.Range("D16").Value = '''2.585.75'
$ws.Range("E16").Value = '''  -1.58%  '
$ws.Range("E17").Value = '''  -3.24%  '
$ws.Range("E18").Value = '''  -0.84%  '
$ws.Range("D19").Value = '''333.73'
$ws.Range("E19").Value = '''  -2.40%  '
$ws.Range("D20").Value = '''9.97'
$ws.Range("E20").Value = '''  -2.39%  '
$ws.Range("D21").Value = '''6.13'
$ws.Range("E21").Value = '''  -4.28%  '
$ws.Range("E22").Value = '''  -0.02%  '
$ws.Range("D23").Value = '''66.47'
$ws.Range("E23").Value = '''  -1.55%  '
$ws.Range("E24").Value = '''  +1.11%  '
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '''  +0.08%  '
$ws.Range("E26").Value = '''  -5.30%  '
$ws.Range("D27").Value = '''6.99'
$ws.Range("E27").Value = '''  -3.54%  '
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '''  +0.06%  '
$ws.Range("D29").Value = '''0.0₃0721'
$ws.Range("E29").Value = '''  -4.30%  '
$ws.Range("D30").Value = '''1.62'
$ws.Range("E30").Value = '''  -2.56%  '
$ws.Range("D31").Value = '''154.97'
$ws.Range("E31").Value = '''  +2.77%  '
$ws.Range("D32").Value = '''5.83'
$ws.Range("E32").Value = '''  -0.82%  '
$ws.Range("E33").Value = '''  -0.79%  '
$ws.Range("E34").Value = '''  -3.92%  '
$ws.Range("D35").Value = '''36.88'
$ws.Range("E35").Value = '''  -1.52%  '
$ws.Range("E36").Value = '''  -4.73%  '
$ws.Range("D37").Value = '''0.824'
$ws.Range("E37").Value = '''  +0.19%  '
$ws.Range("D38").Value = '''0.811'
$ws.Range("E38").Value = '''  -3.18%  '
$ws.Range("E39").Value = '''  -3.85%  '
$ws.Range("D40").Value = '''3.56'
$ws.Range("E40").Value = '''  -0.44%  '
$ws.Range("D41").Value = '''281.36'
$ws.Range("E42").Value = '''  +0.14%  '
$ws.Range("D43").Value = '''10.64'
$ws.Range("E43").Value = '''  -0.84%  '
$ws.Range("D44").Value = '''0.587'
$ws.Range("E44").Value = '''  -2.57%  '
$ws.Range("D45").Value = '''0.0943'
$ws.Range("E45").Value = '''  -0.88%  '
$ws.Range("E46").Value = '''  -0.38%  '
$ws.Range("D47").Value = '''18.17'
$ws.Range("E47").Value = '''  -4.94%  '
$ws.Range("D48").Value = '''0.0224'
$ws.Range("E48").Value = '''  -0.17%  '
$ws.Range("D49").Value = '''1.900.64'
$ws.Range("E49").Value = '''  -3.39%  '
$ws.Range("D50").Value = '''17.65'
$ws.Range("E50").Value = '''  -4.59%  '
$ws.Range("D51").Value = '''4.33'
$ws.Range("E51").Value = '''  -4.93%  '
